$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (3, 4, 5) are being reordered: the non-shared values in
# columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg) rotate between the
# three rows while the rest of each row's content stays the same.

$ws.Range("D3").Value = 44253
$ws.Range("M3").Value = 160

$ws.Range("D4").Value = 44252
$ws.Range("M4").Value = 120
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 13500
$ws.Range("S4").Value = 750

$ws.Range("D5").Value = 44257
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("S5").Value = 806
